$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename metadata4Ing -> metadata4ing, add metadata4ing_DEF column ---
$ws.Range("D1").Value = 'metadata4ing_IRI'
$ws.Range("E1").Value = 'metadata4ing_DESC'
$ws.Range("F1").Value = 'metadata4ing_DEF'
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# --- New column F (metadata4ing_DEF) values for existing rows 2-5 ---
$ws.Range("F2").Value = '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]'
$ws.Range("F3").Value = '[''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]'', ''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])'']'
$ws.Range("F4").Value = '[locstr(''A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.'', ''en'')]'
$ws.Range("F5").Value = '[locstr(''A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.'', ''en'')]'

# --- New row 6: ENVO_03501196 / foaf Group ---
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'http://purl.obolibrary.org/obo/ENVO_03501196'
$ws.Range("C6").Value = '{''label'': ''group'', ''prefLabel'': None, ''altLabel'': None, ''name'': ''ENVO_03501196''}'
$ws.Range("D6").Value = 'http://xmlns.com/foaf/0.1/Group'
$ws.Range("E6").Value = '{''label'': ''group'', ''prefLabel'': ''group'', ''name'': ''group''}'
$ws.Range("F6").Value = '[]'

# --- Hyperlinks for new row 6 (mirrors B2:B5 / D2:D5 pattern) ---
$ws.Hyperlinks.Add($ws.Range("B6"), 'http://purl.obolibrary.org/obo/ENVO_03501196')
$ws.Hyperlinks.Add($ws.Range("D6"), 'http://xmlns.com/foaf/0.1/Group')
# Hyperlinks.Add mints its own style slot; restore the shared Hyperlink style
# so B6/D6 reuse the same cell format (s=2) as the other link cells in the column.
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("D6").Style = "Hyperlink"

Write-Host "done"
